$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.770.77"
$ws.Range("E2").Value = "  +0.03%  "

$ws.Range("D3").Value = "3.560.32"
$ws.Range("E3").Value = "  -0.71%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.68%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "188.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.89%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.629"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.45%  "

$ws.Range("D8").Value = "3.552.57"
$ws.Range("E8").Value = "  -0.82%  "

$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.176"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.655"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.11%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.85%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000299"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.84%  "

$ws.Range("D15").Value = "4.146.43"
$ws.Range("E15").Value = "  -0.23%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.81%  "

$ws.Range("D17").Value = "3.578.46"
$ws.Range("E17").Value = "  -0.19%  "

$ws.Range("D18").Value = "69.714.85"
$ws.Range("E18").Value = "  +0.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.91%  "

$ws.Range("E20").Value = "  +0.42%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "469.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +14.63%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "88.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.22%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.73%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.83%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.03%  "

$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.08%  "

$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.119"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.64%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "65.54"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "583.99"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "38.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.41%  "

$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.13%  "

$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0801"
$ws.Range("E38").Value = "  -4.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.394"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.24%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.51"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.75%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.139"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.56%  "

$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +13.60%  "

$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.32%  "

$ws.Range("D44").Value = "3.212.69"
$ws.Range("E44").Value = "  -3.90%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0440"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.17%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.68%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.136"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.06%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.44%  "
